# TissueWorking.pptx maintenance pass: drop the obsolete PreProcess/Process/
# Conditions/PostProcess flowchart slide (it was the first slide in the deck)
# as part of the documentation cleanup ahead of release.
$p = $ppt.ActivePresentation

$p.Slides.Item(1).Delete()
